$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts existing rows 9-25 down to 10-26)
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with a new weekly data point
# (same market/category/variety as the former row 9, new date 2021-10-26)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = [DateTime]"2021-10-26"
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112013
$ws.Cells.Item(9, 7).Value = "Alcachofa"
$ws.Cells.Item(9, 8).Value = "Madrigal"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 11000
$ws.Cells.Item(9, 12).Value = 12000
$ws.Cells.Item(9, 13).Value = 11500
$ws.Cells.Item(9, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 288
$ws.Cells.Item(9, 17).Value = 40
$ws.Cells.Item(9, 18).Value = "Hortaliza"
